$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Sheet1 currency row (row 6): label lowercase, so new "currency" string is added first
$ws1.Range("A6").Value = "currency"

# Fix product name hyphen on both sheets (adds the hyphenated string next)
$ws1.Range("B1").Value = "821-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-Late Repayment"
$ws2.Range("B1").Value = "821-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-Late Repayment"

# Sheet1 currency value trimmed (adds the trimmed "US Dollar" string last), restyle B6
$ws1.Range("B6").Value = "US Dollar"
$ws1.Range("B6").Interior.Color = 5296274

# Selections (active cell)
$ws2.Range("B1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A6:B6").Select() | Out-Null
